$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: politeness_score becomes a real number (2) instead of text,
# and polite_expressions becomes an empty string instead of "nan".
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = ""

# Row 14: new annotation row appended below row 13.
$ws.Range("A14").Value = "parisk"
$ws.Range("B14").Value = "2"
$ws.Range("C14").Value = "nan"
$ws.Range("D14").Value = "DIS"
$ws.Range("E14").Value = "RES"
$ws.Range("F14").Value = "3a6bf25f-9f71-48b7-a40b-7e968e5f9337"
$ws.Range("G14").Value = "ry-TW-WAb_annotated.xlsx"
$ws.Range("H14").Value = "I suggest to change it to e.g. 'from the true to the approximate posterior' to avoid confusion."
